# Adds a new "Hot Water Fixtures" measure block to the Variables sheet and a
# matching summary row to the Outputs sheet, mirroring the rest of the
# ReportingMeasure blocks already present in the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Outputs sheet: insert a new row at row 40 for the Hot Water Fixtures
#    output, pushing every subsequent row down by one. Written first so the
#    new shared-string entries land in the same order as the authored file.
# ---------------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("Outputs")

$wsOut.Rows.Item(40).Insert(-4121)
$wsOut.Range("A41:I41").Copy()
$wsOut.Range("A40:I40").PasteSpecial(-4122)
$wsOut.Application.CutCopyMode = 0

$wsOut.Range("A40").Value = "Hot Water Fixtures"
$wsOut.Range("D40").Value = "res_stock_reporting.hot_water_fixtures"
$wsOut.Range("F40").Value = "String"
$wsOut.Range("G40").Value = $false
$wsOut.Range("H40").Value = $false
$wsOut.Range("I40").Value = $false

$wsOut.Range("D41").Select()

# ---------------------------------------------------------------------------
# 2) Variables sheet: insert a new 3-row measure block at row 117 (pushing
#    every block that follows down by 3 rows), formatted like its neighbours.
# ---------------------------------------------------------------------------
$wsVar = $wb.Worksheets.Item("Variables")

# Insert 3 blank rows, pulling formatting from the rows below (which, after
# the insert, are the original row-117 block shifted down to 120-122) so the
# new rows inherit the same per-cell styles used by every other block.
$wsVar.Rows.Item(117).Resize(3).Insert(-4121)
$wsVar.Range("A120:Z122").Copy()
$wsVar.Range("A117:Z119").PasteSpecial(-4122)
$wsVar.Application.CutCopyMode = 0

# Row 117 - measure header
$wsVar.Range("A117").Value = $true
$wsVar.Range("B117").Value = "Set Hot Water Fixtures"
$wsVar.Range("C117").Value = "CallMetaMeasure"
$wsVar.Range("D117").Value = "CallMetaMeasure"
$wsVar.Range("E117").Value = "RubyMeasure"

# Row 118 - argument / sample value file
$wsVar.Range("B118").Value = "argument"
$wsVar.Range("D118").Value = "Probability Distributions File"
$wsVar.Range("E118").Value = "probability_file"
$wsVar.Range("G118").Value = "string"

# Row 119 - sample value distribution
$wsVar.Range("B119").Value = "variable"
$wsVar.Range("D119").Value = "Hot Water Fixtures Sample Value"
$wsVar.Range("E119").Value = "sample_value"
$wsVar.Range("G119").Value = "double"
$wsVar.Range("I119").Value = 0.5
$wsVar.Range("K119").Value = 0
$wsVar.Range("L119").Value = 1
$wsVar.Range("M119").Value = 0.5
$wsVar.Range("N119").Value = 0.1666667
$wsVar.Range("R119").Value = "uniform"

# I118 written after D119 so the new shared-string table entries land in the
# same order the workbook was originally authored in.
$wsVar.Range("I118").Value = "hot_water_fixtures.txt"

# Refresh the sheet view to match the authored state.
$wsVar.Application.ActiveWindow.ScrollRow = 93
$wsVar.Range("I119").Select()
